$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 7; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        if ($r -eq 2 -and $c -eq 2) {
            continue
        }
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = "Row $r, Col $c"
    }
}
